$p = $ppt.ActivePresentation
$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = "5/23/25"
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "5/23/25"
        }
    }
}

$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    $sh = $nm.Shapes.Item($j)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = "5/23/25"
    }
}
Write-Host "done"
